$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1281.0952
$ws.Range("I19").Value = 1050.7273
$ws.Range("K19").Value = 1050.7273
$ws.Range("M19").Value = -875.7273
$ws.Range("H21").Value = 19999
$ws.Range("I21").Value = 19999
$ws.Range("K21").Value = 19999
$ws.Range("M21").Value = -19531
$ws.Range("H23").Value = 19999
$ws.Range("I23").Value = 19999
$ws.Range("K23").Value = 19999
$ws.Range("M23").Value = -19765
$ws.Range("H34").Value = 4512.8335
$ws.Range("I34").Value = 4512.8335
$ws.Range("K34").Value = 4512.8335
$ws.Range("M34").Value = -4309.8335
$ws.Range("H36").Value = 4512.8335
$ws.Range("I36").Value = 4512.8335
$ws.Range("K36").Value = 4512.8335
$ws.Range("M36").Value = -3797.8335
$ws.Range("H40").Value = 2679.3076
$ws.Range("J40").Value = 2499.6667
$ws.Range("L40").Value = 2499.6667
$ws.Range("N40").Value = -2849.6667
$ws.Range("H51").Value = 9568.857
$ws.Range("J51").Value = 9561
$ws.Range("L51").Value = 9561
$ws.Range("N51").Value = -10529
$ws.Range("H92").Value = 3980.5833
$ws.Range("I92").Value = 1588.4
$ws.Range("J92").Value = 5689.2856
$ws.Range("K92").Value = 1588.4
$ws.Range("L92").Value = 5689.2856
$ws.Range("M92").Value = -340.4000000000001
$ws.Range("N92").Value = -8185.2856
$ws.Range("H98").Value = 2241.5334
$ws.Range("I98").Value = 1552.1666
$ws.Range("K98").Value = 1552.1666
$ws.Range("M98").Value = -54.16660000000002
$ws.Range("H100").Value = 801.6667
$ws.Range("I100").Value = 801.6667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 801.6667
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -260.6667
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 6019.0557
$ws.Range("I113").Value = 4852.5
$ws.Range("J113").Value = 6952.3
$ws.Range("K113").Value = 4852.5
$ws.Range("L113").Value = 6952.3
$ws.Range("M113").Value = -1598.5
$ws.Range("N113").Value = -13460.3
$ws.Range("H122").Value = 2241.5334
$ws.Range("I122").Value = 1552.1666
$ws.Range("K122").Value = 4656.4998
$ws.Range("M122").Value = -2206.4998
$ws.Range("H127").Value = 2468.4
$ws.Range("I127").Value = 2166.25
$ws.Range("K127").Value = 6498.75
$ws.Range("M127").Value = -1538.75
$ws.Range("H129").Value = 4399
$ws.Range("I129").Value = 5548
$ws.Range("K129").Value = 16644
$ws.Range("M129").Value = -11644
$ws.Range("H131").Value = 3133.3333
$ws.Range("I131").Value = 3133.3333
$ws.Range("K131").Value = 9399.999899999999
$ws.Range("M131").Value = -4359.999899999999
$ws.Range("H132").Value = 402052.2
$ws.Range("I132").Value = 2190.238
$ws.Range("K132").Value = 6570.714
$ws.Range("M132").Value = -4040.714
$ws.Range("H137").Value = 1410.6
$ws.Range("I137").Value = 810.125
$ws.Range("K137").Value = 2430.375
$ws.Range("M137").Value = 119.625
$ws.Range("H138").Value = 2725.8767
$ws.Range("I138").Value = 5533.909
$ws.Range("K138").Value = 16601.727
$ws.Range("M138").Value = -11461.727
$ws.Range("H141").Value = 4361.778
$ws.Range("I141").Value = 2649.077
$ws.Range("K141").Value = 7947.231000000001
$ws.Range("M141").Value = -2767.231000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1607.8
$ws.Range("I2").Value = 1607.8
$ws.Range("K2").Value = 1607.8
$ws.Range("M2").Value = -1494.8
$ws.Range("H5").Value = 1375
$ws.Range("I5").Value = 833.3333
$ws.Range("K5").Value = 833.3333
$ws.Range("M5").Value = -721.3333
$ws.Range("H25").Value = 6666.6665
$ws.Range("I25").Value = 7000
$ws.Range("K25").Value = 7000
$ws.Range("M25").Value = -6598
$ws.Range("H32").Value = 4833.067
$ws.Range("I32").Value = 1730.4615
$ws.Range("K32").Value = 1730.4615
$ws.Range("M32").Value = -1443.4615
$ws.Range("H61").Value = 2876.6667
$ws.Range("I61").Value = 2337.45
$ws.Range("K61").Value = 2337.45
$ws.Range("M61").Value = -2125.45
$ws.Range("H74").Value = 1635.8948
$ws.Range("I74").Value = 1212.7241
$ws.Range("K74").Value = 1212.7241
$ws.Range("M74").Value = -338.7240999999999
$ws.Range("H77").Value = 1635.8948
$ws.Range("I77").Value = 1212.7241
$ws.Range("K77").Value = 6063.620499999999
$ws.Range("M77").Value = -1695.620499999999
$ws.Range("H98").Value = 39999
$ws.Range("J98").Value = 39999
$ws.Range("L98").Value = 39999
$ws.Range("N98").Value = -45989
$ws.Range("H110").Value = 1366.3529
$ws.Range("I110").Value = 1250.2142
$ws.Range("J110").Value = 1908.3334
$ws.Range("K110").Value = 1250.2142
$ws.Range("L110").Value = 1908.3334
$ws.Range("M110").Value = 794.7858000000001
$ws.Range("N110").Value = -5998.3334
$ws.Range("H116").Value = 1607.8
$ws.Range("I116").Value = 1607.8
$ws.Range("K116").Value = 1607.8
$ws.Range("M116").Value = 686.2
$ws.Range("H122").Value = 1615.3846
$ws.Range("I122").Value = 1416.7084
$ws.Range("K122").Value = 4250.1252
$ws.Range("M122").Value = -1800.1252
$ws.Range("H132").Value = 2158.261
$ws.Range("I132").Value = 2179.1765
$ws.Range("K132").Value = 6537.529500000001
$ws.Range("M132").Value = -4007.529500000001
$ws.Range("H136").Value = 2876.6667
$ws.Range("I136").Value = 2337.45
$ws.Range("K136").Value = 7012.349999999999
$ws.Range("M136").Value = -4462.349999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1607.8
$ws.Range("I3").Value = 1607.8
$ws.Range("K3").Value = 1607.8
$ws.Range("M3").Value = -1493.8
$ws.Range("H4").Value = 1375
$ws.Range("I4").Value = 833.3333
$ws.Range("K4").Value = 833.3333
$ws.Range("M4").Value = -718.3333
$ws.Range("H22").Value = 375
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -260.33334
$ws.Range("N22").Value = -546
$ws.Range("H86").Value = 3315.8333
$ws.Range("I86").Value = 1632.6666
$ws.Range("J86").Value = 4999
$ws.Range("K86").Value = 1632.6666
$ws.Range("L86").Value = 4999
$ws.Range("M86").Value = -509.6666
$ws.Range("N86").Value = -7245
$ws.Range("H89").Value = 3315.8333
$ws.Range("I89").Value = 1632.6666
$ws.Range("J89").Value = 4999
$ws.Range("K89").Value = 8163.333000000001
$ws.Range("L89").Value = 24995
$ws.Range("M89").Value = -2547.333000000001
$ws.Range("N89").Value = -36227
$ws.Range("H99").Value = 4036.375
$ws.Range("I99").Value = 4056
$ws.Range("K99").Value = 4056
$ws.Range("M99").Value = -2558
$ws.Range("H105").Value = 2077.889
$ws.Range("I105").Value = 1788.8667
$ws.Range("J105").Value = 2439.1667
$ws.Range("K105").Value = 1788.8667
$ws.Range("L105").Value = 2439.1667
$ws.Range("M105").Value = -41.86670000000004
$ws.Range("N105").Value = -5933.1667
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 635.8
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 594.75
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 594.75
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -1294.75
$ws.Range("H31").Value = 2504.6316
$ws.Range("I31").Value = 1792.2858
$ws.Range("J31").Value = 2920.1667
$ws.Range("K31").Value = 1792.2858
$ws.Range("L31").Value = 2920.1667
$ws.Range("M31").Value = -1497.2858
$ws.Range("N31").Value = -3510.1667
$ws.Range("H34").Value = 2504.6316
$ws.Range("I34").Value = 1792.2858
$ws.Range("J34").Value = 2920.1667
$ws.Range("K34").Value = 1792.2858
$ws.Range("L34").Value = 2920.1667
$ws.Range("M34").Value = -1590.2858
$ws.Range("N34").Value = -3324.1667
$ws.Range("H58").Value = 2870.9285
$ws.Range("I58").Value = 2870.9285
$ws.Range("K58").Value = 2870.9285
$ws.Range("M58").Value = -2667.9285
$ws.Range("H59").Value = 37295.4
$ws.Range("I59").Value = 30749.75
$ws.Range("K59").Value = 30749.75
$ws.Range("M59").Value = -29604.75
$ws.Range("H99").Value = 62499.668
$ws.Range("I99").Value = 8749.5
$ws.Range("K99").Value = 8749.5
$ws.Range("M99").Value = -7251.5
$ws.Range("H107").Value = 1453.6666
$ws.Range("I107").Value = 1282.9375
$ws.Range("K107").Value = 1282.9375
$ws.Range("M107").Value = 637.0625
$ws.Range("H126").Value = 62499.668
$ws.Range("I126").Value = 8749.5
$ws.Range("K126").Value = 26248.5
$ws.Range("M126").Value = -23778.5
$ws.Range("H134").Value = 3039.4736
$ws.Range("I134").Value = 3044.1177
$ws.Range("K134").Value = 9132.3531
$ws.Range("M134").Value = -6597.3531
$ws.Range("H136").Value = 2870.9285
$ws.Range("I136").Value = 2870.9285
$ws.Range("K136").Value = 8612.7855
$ws.Range("M136").Value = -6062.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 253.32
$ws.Range("I12").Value = 318.46155
$ws.Range("K12").Value = 955.38465
$ws.Range("M12").Value = -782.38465
$ws.Range("H14").Value = 174.64285
$ws.Range("I14").Value = 174.64285
$ws.Range("K14").Value = 523.9285500000001
$ws.Range("M14").Value = -350.9285500000001
$ws.Range("H23").Value = 358.66666
$ws.Range("J23").Value = 400
$ws.Range("L23").Value = 1200
$ws.Range("N23").Value = -1670
$ws.Range("H32").Value = 999
$ws.Range("I32").Value = 999
$ws.Range("K32").Value = 2997
$ws.Range("M32").Value = -2714
$ws.Range("H60").Value = 314.66666
$ws.Range("I60").Value = 291.625
$ws.Range("K60").Value = 874.875
$ws.Range("M60").Value = -623.875
$ws.Range("H129").Value = 6149.625
$ws.Range("J129").Value = 6149.625
$ws.Range("L129").Value = 18448.875
$ws.Range("N129").Value = -28448.875
$ws.Range("H131").Value = 21001.828
$ws.Range("J131").Value = 1894.4546
$ws.Range("L131").Value = 5683.3638
$ws.Range("N131").Value = -15763.3638
$ws.Range("H132").Value = 1123.125
$ws.Range("I132").Value = 1121.75
$ws.Range("J132").Value = 1124.5
$ws.Range("K132").Value = 10095.75
$ws.Range("L132").Value = 10120.5
$ws.Range("M132").Value = -7565.75
$ws.Range("N132").Value = -15180.5
$ws.Range("H139").Value = 14583.875
$ws.Range("I139").Value = 9223.25
$ws.Range("J139").Value = 19944.5
$ws.Range("K139").Value = 27669.75
$ws.Range("L139").Value = 59833.5
$ws.Range("M139").Value = -22529.75
$ws.Range("N139").Value = -70113.5
$ws.Range("H140").Value = 8394.5
$ws.Range("I140").Value = 4256.5
$ws.Range("J140").Value = 12532.5
$ws.Range("K140").Value = 12769.5
$ws.Range("L140").Value = 37597.5
$ws.Range("M140").Value = -7589.5
$ws.Range("N140").Value = -47957.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 286.77777
$ws.Range("I2").Value = 260.125
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 260.125
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -147.125
$ws.Range("N2").Value = -726
$ws.Range("H10").Value = 50004.5
$ws.Range("I10").Value = 10
$ws.Range("K10").Value = 10
$ws.Range("M10").Value = 159
$ws.Range("H122").Value = 1944.6364
$ws.Range("I122").Value = 1601.9412
$ws.Range("J122").Value = 3109.8
$ws.Range("K122").Value = 4805.8236
$ws.Range("L122").Value = 9329.400000000001
$ws.Range("M122").Value = -2355.8236
$ws.Range("N122").Value = -14229.4
$ws.Range("H132").Value = 2353.6191
$ws.Range("I132").Value = 2348.6
$ws.Range("K132").Value = 7045.799999999999
$ws.Range("M132").Value = -4515.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 416.72726
$ws.Range("I16").Value = 423.5
$ws.Range("J16").Value = 349
$ws.Range("K16").Value = 423.5
$ws.Range("L16").Value = 349
$ws.Range("M16").Value = -253.5
$ws.Range("N16").Value = -689
$ws.Range("H20").Value = 8441.444
$ws.Range("J20").Value = 8710.429
$ws.Range("L20").Value = 8710.429
$ws.Range("N20").Value = -9162.429
$ws.Range("H22").Value = 1773.25
$ws.Range("I22").Value = 1549
$ws.Range("J22").Value = 1997.5
$ws.Range("K22").Value = 1549
$ws.Range("L22").Value = 1997.5
$ws.Range("M22").Value = -1254
$ws.Range("N22").Value = -2587.5
$ws.Range("H25").Value = 432.66666
$ws.Range("I25").Value = 400
$ws.Range("K25").Value = 400
$ws.Range("M25").Value = -170
$ws.Range("H27").Value = 1773.25
$ws.Range("I27").Value = 1549
$ws.Range("J27").Value = 1997.5
$ws.Range("K27").Value = 1549
$ws.Range("L27").Value = 1997.5
$ws.Range("M27").Value = -1442
$ws.Range("N27").Value = -2211.5
$ws.Range("H42").Value = 20667.5
$ws.Range("I42").Value = 12801
$ws.Range("J42").Value = 60000
$ws.Range("K42").Value = 12801
$ws.Range("L42").Value = 60000
$ws.Range("M42").Value = -12238
$ws.Range("N42").Value = -61126
$ws.Range("H43").Value = 1765976.2
$ws.Range("I43").Value = 30000
$ws.Range("J43").Value = 5733922
$ws.Range("K43").Value = 30000
$ws.Range("L43").Value = 5733922
$ws.Range("M43").Value = -29807
$ws.Range("N43").Value = -5734308
$ws.Range("H46").Value = 3068.4375
$ws.Range("J46").Value = 3499.8462
$ws.Range("L46").Value = 3499.8462
$ws.Range("N46").Value = -3875.8462
$ws.Range("H49").Value = 20667.5
$ws.Range("I49").Value = 12801
$ws.Range("J49").Value = 60000
$ws.Range("K49").Value = 12801
$ws.Range("L49").Value = 60000
$ws.Range("M49").Value = -12654
$ws.Range("N49").Value = -60294
$ws.Range("H122").Value = 3963.5217
$ws.Range("I122").Value = 4607.3076
$ws.Range("K122").Value = 13821.9228
$ws.Range("M122").Value = -11371.9228
$ws.Range("H125").Value = 58333.332
$ws.Range("J125").Value = 58333.332
$ws.Range("L125").Value = 58333.332
$ws.Range("N125").Value = -68173.33199999999
$ws.Range("H132").Value = 2848.0833
$ws.Range("I132").Value = 2316.2222
$ws.Range("K132").Value = 6948.6666
$ws.Range("M132").Value = -4418.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 512916.66
$ws.Range("I3").Value = 1007166.7
$ws.Range("K3").Value = 1007166.7
$ws.Range("M3").Value = -1007052.7
$ws.Range("H4").Value = 19003722
$ws.Range("I4").Value = 19003722
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 19003722
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -19003609
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 25000
$ws.Range("I9").Value = 25000
$ws.Range("K9").Value = 25000
$ws.Range("M9").Value = -24860
$ws.Range("H10").Value = 1006
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1006
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1006
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -1344
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10284
$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19765
$ws.Range("H22").Value = 8000
$ws.Range("J22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("N22").Value = -8586
$ws.Range("H29").Value = 49000
$ws.Range("I29").Value = 49000
$ws.Range("K29").Value = 49000
$ws.Range("M29").Value = -48710
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 30307.834
$ws.Range("I32").Value = 38961.75
$ws.Range("J32").Value = 13000
$ws.Range("K32").Value = 38961.75
$ws.Range("L32").Value = 13000
$ws.Range("M32").Value = -38644.75
$ws.Range("N32").Value = -13634
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 45937.25
$ws.Range("I34").Value = 47874.5
$ws.Range("J34").Value = 44000
$ws.Range("K34").Value = 47874.5
$ws.Range("L34").Value = 44000
$ws.Range("M34").Value = -47671.5
$ws.Range("N34").Value = -44406
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19710
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H37").Value = 53333.332
$ws.Range("I37").Value = 50000
$ws.Range("J37").Value = 60000
$ws.Range("K37").Value = 50000
$ws.Range("L37").Value = 60000
$ws.Range("M37").Value = -49797
$ws.Range("N37").Value = -60406
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H42").Value = 39998
$ws.Range("I42").Value = 39998
$ws.Range("K42").Value = 39998
$ws.Range("M42").Value = -39620
$ws.Range("H43").Value = 42999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 42999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 42999.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -43297.5
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H54").Value = 29250
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 29250
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 29250
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -30290
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 686.3077
$ws.Range("I107").Value = 464.2857
$ws.Range("K107").Value = 1392.8571
$ws.Range("M107").Value = 527.1428999999998
$ws.Range("H122").Value = 5742.3057
$ws.Range("I122").Value = 6555.4346
$ws.Range("J122").Value = 4303.6924
$ws.Range("K122").Value = 19666.3038
$ws.Range("L122").Value = 12911.0772
$ws.Range("M122").Value = -17216.3038
$ws.Range("N122").Value = -17811.0772
$ws.Range("H123").Value = 67500
$ws.Range("J123").Value = 67500
$ws.Range("L123").Value = 67500
$ws.Range("N123").Value = -77300
$ws.Range("H125").Value = 43455
$ws.Range("J125").Value = 43455
$ws.Range("L125").Value = 43455
$ws.Range("N125").Value = -53295
$ws.Range("H136").Value = 3113.2163
$ws.Range("I136").Value = 2152.6177
$ws.Range("K136").Value = 6457.853099999999
$ws.Range("M136").Value = -3907.853099999999
